# Updated cryptos list on Sun May 21 18:37:16 UTC 2023 with GitHub Actions
#
# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the crypto table with
# the latest scrape. Both columns hold pre-formatted strings (D has locale-style
# thousand separators / fixed decimals, E has padded "  +x.xx%  " text), so each
# D write is briefly forced to Text number-format to stop Excel from reinterpreting
# it as a number (which would e.g. collapse "1.010" -> "1.01"); the format flag is
# cleared right after so the cell keeps the workbook's original (unstyled) look.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.093.42'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.05%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.823.36'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.18%  '

$ws.Range("E4").Value = '  -0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.03'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.010'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4623'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -2.33%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3643'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.58%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07299'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.16%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8699'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.83%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.12'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.89%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.867.89'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.37%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07592'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +3.15%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.347'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.57%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.32'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.34%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.488'
$ws.Range("D16").ClearFormats()

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.009'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.39%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008627'
$ws.Range("D18").ClearFormats()

$ws.Range("E19").Value = '  -0.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.414.39'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.03%  '

$ws.Range("E21").Value = '  -2.76%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.207'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.39%  '

$ws.Range("E23").Value = '  -1.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.094.31'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.25%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.96'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.70%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.863'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.29%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.25'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.38%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.091'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -4.58%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.082'
$ws.Range("D30").ClearFormats()

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08905'
$ws.Range("D31").ClearFormats()

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.962'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.54%  '

$ws.Range("E33").Value = '  -4.01%  '

$ws.Range("E34").Value = '  -2.61%  '

$ws.Range("E35").Value = '  -3.63%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.011'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.10%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.477'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.94%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.074'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.18%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05251'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.05%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01915'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.58%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.932'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.44%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.136'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.21%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5197'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.32%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1631'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.36%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.264'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -3.89%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4866'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.17%  '

$ws.Range("E47").Value = '  -0.19%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.18'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.25%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '103.61'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.94%  '

$ws.Range("E50").Value = '  -3.13%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06261'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.04%  '
